# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
# The sheet only tracked per-player stats before; this backfills the
# team's overall season record (98-64-0) across every player row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells, one column past the existing "Unnamed: 28" (AC).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Give the new headers the same bold/centered/bordered look as the rest
# of row 1 by copying the formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill in the team's season record (98 wins, 64 losses, 0 ties) for
# every player row.
$lastRow = 45
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 98
    $ws.Cells.Item($r, 31).Value = 64
    $ws.Cells.Item($r, 32).Value = 0
}
